$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new structure "EnzymeKinematic<T>" is documented on the interface
# coverage sheet. Insert its row right above "ErrorMargin<T>" (current
# row 5), pushing every following row down by one.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "EnzymeKinematic<T>"
$ws.Range("B5").Value = "Oui"
$ws.Range("C5").Value = "S/O"
$ws.Range("D5").Value = "Oui"
$ws.Range("E5").Value = "S/O"
$ws.Range("F5").Value = "Oui"
$ws.Range("G5").Value = "S/O"
$ws.Range("H5").Value = "S/O"
$ws.Range("I5").Value = "S/O"
$ws.Range("J5").Value = "S/O"
$ws.Range("K5").Value = "S/O"
$ws.Range("L5").Value = "Oui"

# Leave the selection where the author ended up after entering the data.
$ws.Range("J17").Select() | Out-Null
